$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Discord Name values: add underscore before the '#'
$ws.Range("D2").Value = "pm_me_cute_sloths_#5223"
$ws.Range("D3").Value = "pm_me_cute_sloths_#5223"

# Move the selection to D4, preparing for the next entry (kickoff)
$ws.Range("D4").Select()
